$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert a new worksheet right after Feuil1 and name it Feuil2
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Feuil2"

# Populate the header row with the new shared strings
# (write in this order so the shared-string table is built as
#  0: capacité, 1: recharge, 2: autonomie)
$ws2.Range("A1").Value = "capacité"
$ws2.Range("C1").Value = "recharge"
$ws2.Range("B1").Value = "autonomie"

# Select cell B1 and make Feuil2 the active sheet/tab
$ws2.Range("B1").Select()
$ws2.Activate()
